$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the latest GitHub Actions refresh.
# Cells whose new value would otherwise be auto-parsed as a number by Excel
# are forced to Text format first so they round-trip as plain text, matching
# how this sheet stores all Price/Volume columns (inline strings).

$ws.Range('D2').Value = '67.393.06'
$ws.Range('E2').Value = '  +7.31%  '
$ws.Range('D3').Value = '3.585.05'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '415.59'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.22'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  +3.71%  '
$ws.Range('D8').Value = '3.578.55'
$ws.Range('E8').Value = '  +3.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.771'
$ws.Range('E10').Value = '  +6.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.181'
$ws.Range('E11').Value = '  +17.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000337'
$ws.Range('E12').Value = '  +48.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.32'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '4.157.74'
$ws.Range('E15').Value = '  +3.47%  '
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.39'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '3.606.73'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.14'
$ws.Range('E19').Value = '  +5.62%  '
$ws.Range('D20').Value = '67.223.23'
$ws.Range('E20').Value = '  +7.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.25'
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '452.06'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '89.11'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('E24').Value = '  -4.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.14'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.34'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('E27').Value = '  -6.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '35.04'
$ws.Range('E28').Value = '  +5.00%  '
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.77'
$ws.Range('E30').Value = '  +3.78%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.35'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('E32').Value = '  +4.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.38'
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.162'
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.89'
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.77'
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').Value = '0.0₃0775'
$ws.Range('E39').Value = '  +36.13%  '
$ws.Range('E40').Value = '  +9.58%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.06'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '149.14'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.314'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.26'
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.31'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  -4.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '115.83'
$ws.Range('E50').Value = '  +6.51%  '
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.59'
$ws.Range('E51').Value = '  +10.76%  '
